# Sistema de legalizacion - actualizar numeros de factura
$wb = $excel.ActiveWorkbook
$wsHeaders = $wb.Worksheets.Item("Headers")
$wsDetails = $wb.Worksheets.Item("Details")

# Headers sheet: update INVOICE_NUMBER (column A) for the 3 invoices
$wsHeaders.Range("A2").Value = "2025-0115"
$wsHeaders.Range("A3").Value = "2025-0116"
$wsHeaders.Range("A4").Value = "2025-0117"

# Details sheet: update INVOICE_NUMBER (column A) for each detail line,
# keeping them grouped per invoice (3 lines each)
$wsDetails.Range("A2").Value = "2025-0115"
$wsDetails.Range("A3").Value = "2025-0115"
$wsDetails.Range("A4").Value = "2025-0115"
$wsDetails.Range("A5").Value = "2025-0116"
$wsDetails.Range("A6").Value = "2025-0116"
$wsDetails.Range("A7").Value = "2025-0116"
$wsDetails.Range("A8").Value = "2025-0117"
$wsDetails.Range("A9").Value = "2025-0117"
$wsDetails.Range("A10").Value = "2025-0117"

# Leave the selection on the Headers sheet at A4 (last place the user
# clicked before moving on), then switch to and activate the Details
# sheet, which keeps its own prior selection (A8:A10).
$wsHeaders.Range("A4").Select()
$wsDetails.Activate()
